$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-8 from 2023-09-01 (45170)
# to 2023-09-05 (45174), matching the automatic update reflected in the diff.
$ws.Range("C2:C8").Value = 45174
